$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SheetWine"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SheetVineyard"

$ws2.Range("A1").Value = "did"
$ws2.Range("B1").Value = "district"
$ws2.Range("A2").Value = "A"
$ws2.Range("B2").Value = "ブルゴーニュ"
$ws2.Range("A3").Value = "B"
$ws2.Range("B3").Value = "ボルドー"
$ws2.Range("A4").Value = "C"
$ws2.Range("B4").Value = "ロワール"
$ws2.Range("A5").Value = "D"
$ws2.Range("B5").Value = "シャンバーニュ"
$ws2.Range("A6").Value = "E"
$ws2.Range("B6").Value = "チリ"

$ws2.Activate()
$ws2.Range("A1:B6").Select()
